# Rerun and summarise models without urban landuse
# Renames each summary sheet, drops the 'Education[T.Secondary+Matura]' row
# (removed reference-category row) and rewrites the coefficient/p-value
# columns with the new model's results.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: summ25494940 -> summ48038137 ---
$ws = $wb.Worksheets.Item(1)
$ws.Name = "summ48038137"

# Drop the old row 3 (Education[T.Secondary+Matura] reference row removed
# from the new model); this shifts all rows below it up by one and
# updates the used range automatically.
$ws.Rows.Item(3).Delete()

# Write the new model's coefficients/p-values (rows 2-16).
$ws.Range("A2").Value = "Intercept"
$ws.Range("B2").Value = [double]"15564.00642328113"
$ws.Range("C2").Value = [double]"2.040045732885005e-09"
$ws.Range("A3").Value = "Education[T.University]"
$ws.Range("B3").Value = [double]"-391.0937004115779"
$ws.Range("C3").Value = [double]"0.2519503641231707"
$ws.Range("A4").Value = "Education[T.Unknown/Other]"
$ws.Range("B4").Value = [double]"-80.90200067768023"
$ws.Range("C4").Value = [double]"0.8255788398047164"
$ws.Range("A5").Value = "HHSize"
$ws.Range("B5").Value = [double]"-168.2265084120404"
$ws.Range("C5").Value = [double]"0.1938187430840698"
$ws.Range("A6").Value = "Sex"
$ws.Range("B6").Value = [double]"-1617.915941088716"
$ws.Range("C6").Value = [double]"5.321938716298819e-09"
$ws.Range("A7").Value = "Age"
$ws.Range("B7").Value = [double]"2.207544243629645"
$ws.Range("C7").Value = [double]"0.8615698987650022"
$ws.Range("A8").Value = "DistSubcenter_res"
$ws.Range("B8").Value = [double]"-411.5918882146606"
$ws.Range("C8").Value = [double]"1.023688495886929e-05"
$ws.Range("A9").Value = "DistCenter_res"
$ws.Range("B9").Value = [double]"735.6113021815437"
$ws.Range("C9").Value = [double]"2.262765448153656e-19"
$ws.Range("A10").Value = "UrbPopDensity_res"
$ws.Range("B10").Value = [double]"0.08356502631882656"
$ws.Range("C10").Value = [double]"0.1480899081835948"
$ws.Range("A11").Value = "UrbBuildDensity_res"
$ws.Range("B11").Value = [double]"9.110529732999537e-06"
$ws.Range("C11").Value = [double]"0.7086303400114848"
$ws.Range("A12").Value = "IntersecDensity_res"
$ws.Range("B12").Value = [double]"-32.2227939562399"
$ws.Range("C12").Value = [double]"0.03135170956501851"
$ws.Range("A13").Value = "street_length_res"
$ws.Range("B13").Value = [double]"-20.37283180905422"
$ws.Range("C13").Value = [double]"0.1787157302695445"
$ws.Range("A14").Value = "LU_Comm_res"
$ws.Range("B14").Value = [double]"-4297.976554692949"
$ws.Range("C14").Value = [double]"0.04780802650029996"
$ws.Range("A15").Value = "LU_UrbFab_res"
$ws.Range("B15").Value = [double]"-6004.602961493107"
$ws.Range("C15").Value = [double]"0.002204087660514412"
$ws.Range("A16").Value = "bike_lane_share_res"
$ws.Range("B16").Value = [double]"83.96556927551251"
$ws.Range("C16").Value = [double]"0.9770314929306146"

# --- Sheet 2: summ25808391 -> summ48431239 ---
$ws = $wb.Worksheets.Item(2)
$ws.Name = "summ48431239"

# Drop the old row 3 (Education[T.Secondary+Matura] reference row removed
# from the new model); this shifts all rows below it up by one and
# updates the used range automatically.
$ws.Rows.Item(3).Delete()

# Write the new model's coefficients/p-values (rows 2-16).
$ws.Range("A2").Value = "Intercept"
$ws.Range("B2").Value = [double]"15161.69836584974"
$ws.Range("C2").Value = [double]"2.898898039229992e-09"
$ws.Range("A3").Value = "Education[T.University]"
$ws.Range("B3").Value = [double]"-127.6767808276666"
$ws.Range("C3").Value = [double]"0.7055619708329022"
$ws.Range("A4").Value = "Education[T.Unknown/Other]"
$ws.Range("B4").Value = [double]"203.0291286338824"
$ws.Range("C4").Value = [double]"0.5792910208958866"
$ws.Range("A5").Value = "HHSize"
$ws.Range("B5").Value = [double]"-131.7938497435691"
$ws.Range("C5").Value = [double]"0.301959430199405"
$ws.Range("A6").Value = "Sex"
$ws.Range("B6").Value = [double]"-1645.611790983491"
$ws.Range("C6").Value = [double]"2.194742067547877e-09"
$ws.Range("A7").Value = "Age"
$ws.Range("B7").Value = [double]"-5.78934289556878"
$ws.Range("C7").Value = [double]"0.646971635516374"
$ws.Range("A8").Value = "DistSubcenter_res"
$ws.Range("B8").Value = [double]"-451.7076328913233"
$ws.Range("C8").Value = [double]"9.694039469065625e-07"
$ws.Range("A9").Value = "DistCenter_res"
$ws.Range("B9").Value = [double]"777.030895345513"
$ws.Range("C9").Value = [double]"2.821228555979995e-22"
$ws.Range("A10").Value = "UrbPopDensity_res"
$ws.Range("B10").Value = [double]"0.05640526837733201"
$ws.Range("C10").Value = [double]"0.3140336584510198"
$ws.Range("A11").Value = "UrbBuildDensity_res"
$ws.Range("B11").Value = [double]"3.761512228635299e-07"
$ws.Range("C11").Value = [double]"0.9873095317558158"
$ws.Range("A12").Value = "IntersecDensity_res"
$ws.Range("B12").Value = [double]"-24.62325072798235"
$ws.Range("C12").Value = [double]"0.09498140981080634"
$ws.Range("A13").Value = "street_length_res"
$ws.Range("B13").Value = [double]"-21.54187441590123"
$ws.Range("C13").Value = [double]"0.1491279283367027"
$ws.Range("A14").Value = "LU_Comm_res"
$ws.Range("B14").Value = [double]"-3244.287588647827"
$ws.Range("C14").Value = [double]"0.1296531666331054"
$ws.Range("A15").Value = "LU_UrbFab_res"
$ws.Range("B15").Value = [double]"-5963.2801289904"
$ws.Range("C15").Value = [double]"0.002056442357301361"
$ws.Range("A16").Value = "bike_lane_share_res"
$ws.Range("B16").Value = [double]"2118.498197719266"
$ws.Range("C16").Value = [double]"0.4637804926372833"

# --- Sheet 3: summ26116938 -> summ48815714 ---
$ws = $wb.Worksheets.Item(3)
$ws.Name = "summ48815714"

# Drop the old row 3 (Education[T.Secondary+Matura] reference row removed
# from the new model); this shifts all rows below it up by one and
# updates the used range automatically.
$ws.Rows.Item(3).Delete()

# Write the new model's coefficients/p-values (rows 2-16).
$ws.Range("A2").Value = "Intercept"
$ws.Range("B2").Value = [double]"13032.73996673841"
$ws.Range("C2").Value = [double]"3.951090839748804e-07"
$ws.Range("A3").Value = "Education[T.University]"
$ws.Range("B3").Value = [double]"-281.2241643911053"
$ws.Range("C3").Value = [double]"0.4098941131007241"
$ws.Range("A4").Value = "Education[T.Unknown/Other]"
$ws.Range("B4").Value = [double]"-298.3023388996229"
$ws.Range("C4").Value = [double]"0.4099326254338066"
$ws.Range("A5").Value = "HHSize"
$ws.Range("B5").Value = [double]"-122.5124130333592"
$ws.Range("C5").Value = [double]"0.3405219275481347"
$ws.Range("A6").Value = "Sex"
$ws.Range("B6").Value = [double]"-1679.171208983657"
$ws.Range("C6").Value = [double]"8.853024647622108e-10"
$ws.Range("A7").Value = "Age"
$ws.Range("B7").Value = [double]"6.785339239894712"
$ws.Range("C7").Value = [double]"0.5928086242397401"
$ws.Range("A8").Value = "DistSubcenter_res"
$ws.Range("B8").Value = [double]"-444.3896964575188"
$ws.Range("C8").Value = [double]"1.495499499738234e-06"
$ws.Range("A9").Value = "DistCenter_res"
$ws.Range("B9").Value = [double]"777.2957195611036"
$ws.Range("C9").Value = [double]"9.407951691980246e-22"
$ws.Range("A10").Value = "UrbPopDensity_res"
$ws.Range("B10").Value = [double]"0.1019619045290361"
$ws.Range("C10").Value = [double]"0.07251222809770269"
$ws.Range("A11").Value = "UrbBuildDensity_res"
$ws.Range("B11").Value = [double]"-8.5663535695798e-06"
$ws.Range("C11").Value = [double]"0.7213132574352079"
$ws.Range("A12").Value = "IntersecDensity_res"
$ws.Range("B12").Value = [double]"-15.64290359491772"
$ws.Range("C12").Value = [double]"0.2966453739160945"
$ws.Range("A13").Value = "street_length_res"
$ws.Range("B13").Value = [double]"-8.312472656975068"
$ws.Range("C13").Value = [double]"0.5807591798392753"
$ws.Range("A14").Value = "LU_Comm_res"
$ws.Range("B14").Value = [double]"-3264.438132422921"
$ws.Range("C14").Value = [double]"0.1219039088159581"
$ws.Range("A15").Value = "LU_UrbFab_res"
$ws.Range("B15").Value = [double]"-6713.40638422446"
$ws.Range("C15").Value = [double]"0.0004718287684443202"
$ws.Range("A16").Value = "bike_lane_share_res"
$ws.Range("B16").Value = [double]"-38.02030109243424"
$ws.Range("C16").Value = [double]"0.9894165385810936"

# --- Sheet 4: summ26427482 -> summ49194290 ---
$ws = $wb.Worksheets.Item(4)
$ws.Name = "summ49194290"

# Drop the old row 3 (Education[T.Secondary+Matura] reference row removed
# from the new model); this shifts all rows below it up by one and
# updates the used range automatically.
$ws.Rows.Item(3).Delete()

# Write the new model's coefficients/p-values (rows 2-16).
$ws.Range("A2").Value = "Intercept"
$ws.Range("B2").Value = [double]"13214.81636286243"
$ws.Range("C2").Value = [double]"3.716855497206852e-07"
$ws.Range("A3").Value = "Education[T.University]"
$ws.Range("B3").Value = [double]"-21.09210037404648"
$ws.Range("C3").Value = [double]"0.9508280059482381"
$ws.Range("A4").Value = "Education[T.Unknown/Other]"
$ws.Range("B4").Value = [double]"19.20205287426006"
$ws.Range("C4").Value = [double]"0.9584874912515369"
$ws.Range("A5").Value = "HHSize"
$ws.Range("B5").Value = [double]"-68.44637837920331"
$ws.Range("C5").Value = [double]"0.5944103774745833"
$ws.Range("A6").Value = "Sex"
$ws.Range("B6").Value = [double]"-1548.136855769964"
$ws.Range("C6").Value = [double]"1.867852461868582e-08"
$ws.Range("A7").Value = "Age"
$ws.Range("B7").Value = [double]"-5.523145479820599"
$ws.Range("C7").Value = [double]"0.6649028161996674"
$ws.Range("A8").Value = "DistSubcenter_res"
$ws.Range("B8").Value = [double]"-375.705318199687"
$ws.Range("C8").Value = [double]"5.327182483665493e-05"
$ws.Range("A9").Value = "DistCenter_res"
$ws.Range("B9").Value = [double]"809.4932109542676"
$ws.Range("C9").Value = [double]"1.011003710622632e-23"
$ws.Range("A10").Value = "UrbPopDensity_res"
$ws.Range("B10").Value = [double]"0.09054337718822074"
$ws.Range("C10").Value = [double]"0.1147702871106723"
$ws.Range("A11").Value = "UrbBuildDensity_res"
$ws.Range("B11").Value = [double]"-1.12169990053856e-05"
$ws.Range("C11").Value = [double]"0.6315764110588165"
$ws.Range("A12").Value = "IntersecDensity_res"
$ws.Range("B12").Value = [double]"-14.43548996283053"
$ws.Range("C12").Value = [double]"0.3353925729068744"
$ws.Range("A13").Value = "street_length_res"
$ws.Range("B13").Value = [double]"-15.40236575347479"
$ws.Range("C13").Value = [double]"0.3137257893482888"
$ws.Range("A14").Value = "LU_Comm_res"
$ws.Range("B14").Value = [double]"-1835.883723185776"
$ws.Range("C14").Value = [double]"0.3909412403463206"
$ws.Range("A15").Value = "LU_UrbFab_res"
$ws.Range("B15").Value = [double]"-6853.446959914633"
$ws.Range("C15").Value = [double]"0.0003742233526645857"
$ws.Range("A16").Value = "bike_lane_share_res"
$ws.Range("B16").Value = [double]"1185.398098222433"
$ws.Range("C16").Value = [double]"0.6804529694093098"

# --- Sheet 5: summ26747142 -> summ49655699 ---
$ws = $wb.Worksheets.Item(5)
$ws.Name = "summ49655699"

# Drop the old row 3 (Education[T.Secondary+Matura] reference row removed
# from the new model); this shifts all rows below it up by one and
# updates the used range automatically.
$ws.Rows.Item(3).Delete()

# Write the new model's coefficients/p-values (rows 2-16).
$ws.Range("A2").Value = "Intercept"
$ws.Range("B2").Value = [double]"15534.00367424062"
$ws.Range("C2").Value = [double]"2.444823708909066e-09"
$ws.Range("A3").Value = "Education[T.University]"
$ws.Range("B3").Value = [double]"-359.7334155573743"
$ws.Range("C3").Value = [double]"0.297433358913709"
$ws.Range("A4").Value = "Education[T.Unknown/Other]"
$ws.Range("B4").Value = [double]"-393.0175901806799"
$ws.Range("C4").Value = [double]"0.2874578355049398"
$ws.Range("A5").Value = "HHSize"
$ws.Range("B5").Value = [double]"-97.20655465418244"
$ws.Range("C5").Value = [double]"0.4508733138221021"
$ws.Range("A6").Value = "Sex"
$ws.Range("B6").Value = [double]"-1632.107937466568"
$ws.Range("C6").Value = [double]"4.681868904897734e-09"
$ws.Range("A7").Value = "Age"
$ws.Range("B7").Value = [double]"1.913248636913671"
$ws.Range("C7").Value = [double]"0.880985922739152"
$ws.Range("A8").Value = "DistSubcenter_res"
$ws.Range("B8").Value = [double]"-399.491043170475"
$ws.Range("C8").Value = [double]"1.855484426980094e-05"
$ws.Range("A9").Value = "DistCenter_res"
$ws.Range("B9").Value = [double]"729.7272075851233"
$ws.Range("C9").Value = [double]"2.830260560663993e-19"
$ws.Range("A10").Value = "UrbPopDensity_res"
$ws.Range("B10").Value = [double]"0.07126099610638251"
$ws.Range("C10").Value = [double]"0.215883842376783"
$ws.Range("A11").Value = "UrbBuildDensity_res"
$ws.Range("B11").Value = [double]"-4.953738121785727e-06"
$ws.Range("C11").Value = [double]"0.8345254609670845"
$ws.Range("A12").Value = "IntersecDensity_res"
$ws.Range("B12").Value = [double]"-25.12169546440677"
$ws.Range("C12").Value = [double]"0.09394003214957569"
$ws.Range("A13").Value = "street_length_res"
$ws.Range("B13").Value = [double]"-22.17334354942553"
$ws.Range("C13").Value = [double]"0.1468043537315911"
$ws.Range("A14").Value = "LU_Comm_res"
$ws.Range("B14").Value = [double]"-3190.289053143269"
$ws.Range("C14").Value = [double]"0.1428125299445207"
$ws.Range("A15").Value = "LU_UrbFab_res"
$ws.Range("B15").Value = [double]"-6589.956495129978"
$ws.Range("C15").Value = [double]"0.0007476528499656639"
$ws.Range("A16").Value = "bike_lane_share_res"
$ws.Range("B16").Value = [double]"1234.54254164007"
$ws.Range("C16").Value = [double]"0.6746305629752161"

# --- Sheet 6: summ27069650 -> summ50021425 ---
$ws = $wb.Worksheets.Item(6)
$ws.Name = "summ50021425"

# Drop the old row 3 (Education[T.Secondary+Matura] reference row removed
# from the new model); this shifts all rows below it up by one and
# updates the used range automatically.
$ws.Rows.Item(3).Delete()

# Write the new model's coefficients/p-values (rows 2-16).
$ws.Range("A2").Value = "Intercept"
$ws.Range("B2").Value = [double]"15088.35562453749"
$ws.Range("C2").Value = [double]"3.814812522749387e-09"
$ws.Range("A3").Value = "Education[T.University]"
$ws.Range("B3").Value = [double]"-179.3493055126713"
$ws.Range("C3").Value = [double]"0.6005091597027234"
$ws.Range("A4").Value = "Education[T.Unknown/Other]"
$ws.Range("B4").Value = [double]"-159.3001505028262"
$ws.Range("C4").Value = [double]"0.661379191019369"
$ws.Range("A5").Value = "HHSize"
$ws.Range("B5").Value = [double]"-88.50904660506973"
$ws.Range("C5").Value = [double]"0.4871604104883421"
$ws.Range("A6").Value = "Sex"
$ws.Range("B6").Value = [double]"-1797.382374401847"
$ws.Range("C6").Value = [double]"7.421636235179066e-11"
$ws.Range("A7").Value = "Age"
$ws.Range("B7").Value = [double]"0.322667354118459"
$ws.Range("C7").Value = [double]"0.9796863375047218"
$ws.Range("A8").Value = "DistSubcenter_res"
$ws.Range("B8").Value = [double]"-446.0341821235484"
$ws.Range("C8").Value = [double]"1.659786644672771e-06"
$ws.Range("A9").Value = "DistCenter_res"
$ws.Range("B9").Value = [double]"848.6218747620194"
$ws.Range("C9").Value = [double]"1.898895854356808e-25"
$ws.Range("A10").Value = "UrbPopDensity_res"
$ws.Range("B10").Value = [double]"0.1144653348056466"
$ws.Range("C10").Value = [double]"0.04681533452252361"
$ws.Range("A11").Value = "UrbBuildDensity_res"
$ws.Range("B11").Value = [double]"2.170298152243693e-06"
$ws.Range("C11").Value = [double]"0.9256460905461489"
$ws.Range("A12").Value = "IntersecDensity_res"
$ws.Range("B12").Value = [double]"-29.89014352248457"
$ws.Range("C12").Value = [double]"0.04420030506635398"
$ws.Range("A13").Value = "street_length_res"
$ws.Range("B13").Value = [double]"-23.9716735797459"
$ws.Range("C13").Value = [double]"0.1102034499590424"
$ws.Range("A14").Value = "LU_Comm_res"
$ws.Range("B14").Value = [double]"-3696.709361166511"
$ws.Range("C14").Value = [double]"0.08549699236022656"
$ws.Range("A15").Value = "LU_UrbFab_res"
$ws.Range("B15").Value = [double]"-5821.699188706059"
$ws.Range("C15").Value = [double]"0.002885951682545536"
$ws.Range("A16").Value = "bike_lane_share_res"
$ws.Range("B16").Value = [double]"-507.8660807429724"
$ws.Range("C16").Value = [double]"0.861555132809653"

# --- Sheet 7: summ27412015 -> summ50367714 ---
$ws = $wb.Worksheets.Item(7)
$ws.Name = "summ50367714"

# Drop the old row 3 (Education[T.Secondary+Matura] reference row removed
# from the new model); this shifts all rows below it up by one and
# updates the used range automatically.
$ws.Rows.Item(3).Delete()

# Write the new model's coefficients/p-values (rows 2-16).
$ws.Range("A2").Value = "Intercept"
$ws.Range("B2").Value = [double]"15143.12242992946"
$ws.Range("C2").Value = [double]"4.128947200603713e-09"
$ws.Range("A3").Value = "Education[T.University]"
$ws.Range("B3").Value = [double]"-317.3608175456213"
$ws.Range("C3").Value = [double]"0.3594282554236714"
$ws.Range("A4").Value = "Education[T.Unknown/Other]"
$ws.Range("B4").Value = [double]"-330.816668400657"
$ws.Range("C4").Value = [double]"0.3717394724384717"
$ws.Range("A5").Value = "HHSize"
$ws.Range("B5").Value = [double]"-146.4311758038386"
$ws.Range("C5").Value = [double]"0.2553578062775655"
$ws.Range("A6").Value = "Sex"
$ws.Range("B6").Value = [double]"-1546.247534808731"
$ws.Range("C6").Value = [double]"2.247846951675864e-08"
$ws.Range("A7").Value = "Age"
$ws.Range("B7").Value = [double]"-7.51880873455982"
$ws.Range("C7").Value = [double]"0.5544042221149315"
$ws.Range("A8").Value = "DistSubcenter_res"
$ws.Range("B8").Value = [double]"-371.70070532663"
$ws.Range("C8").Value = [double]"7.093918575471537e-05"
$ws.Range("A9").Value = "DistCenter_res"
$ws.Range("B9").Value = [double]"741.0150574015516"
$ws.Range("C9").Value = [double]"8.069350188706726e-20"
$ws.Range("A10").Value = "UrbPopDensity_res"
$ws.Range("B10").Value = [double]"0.05997553256932597"
$ws.Range("C10").Value = [double]"0.2941706097978063"
$ws.Range("A11").Value = "UrbBuildDensity_res"
$ws.Range("B11").Value = [double]"9.45777047893102e-07"
$ws.Range("C11").Value = [double]"0.9683397082204751"
$ws.Range("A12").Value = "IntersecDensity_res"
$ws.Range("B12").Value = [double]"-22.21772885465363"
$ws.Range("C12").Value = [double]"0.1380193868693297"
$ws.Range("A13").Value = "street_length_res"
$ws.Range("B13").Value = [double]"-19.61355064392484"
$ws.Range("C13").Value = [double]"0.1927303607426932"
$ws.Range("A14").Value = "LU_Comm_res"
$ws.Range("B14").Value = [double]"-4225.014943990814"
$ws.Range("C14").Value = [double]"0.04852868639730452"
$ws.Range("A15").Value = "LU_UrbFab_res"
$ws.Range("B15").Value = [double]"-6424.500252837506"
$ws.Range("C15").Value = [double]"0.001114566128041987"
$ws.Range("A16").Value = "bike_lane_share_res"
$ws.Range("B16").Value = [double]"2438.90346190184"
$ws.Range("C16").Value = [double]"0.4055508915084639"

# --- Sheet 8: summ27712805 -> summ50703252 ---
$ws = $wb.Worksheets.Item(8)
$ws.Name = "summ50703252"

# Drop the old row 3 (Education[T.Secondary+Matura] reference row removed
# from the new model); this shifts all rows below it up by one and
# updates the used range automatically.
$ws.Rows.Item(3).Delete()

# Write the new model's coefficients/p-values (rows 2-16).
$ws.Range("A2").Value = "Intercept"
$ws.Range("B2").Value = [double]"15404.47170854902"
$ws.Range("C2").Value = [double]"5.069111938054369e-09"
$ws.Range("A3").Value = "Education[T.University]"
$ws.Range("B3").Value = [double]"-311.5828529286151"
$ws.Range("C3").Value = [double]"0.3723229738774081"
$ws.Range("A4").Value = "Education[T.Unknown/Other]"
$ws.Range("B4").Value = [double]"-108.6400863387226"
$ws.Range("C4").Value = [double]"0.7704694093587219"
$ws.Range("A5").Value = "HHSize"
$ws.Range("B5").Value = [double]"-133.6951344126564"
$ws.Range("C5").Value = [double]"0.309468435943858"
$ws.Range("A6").Value = "Sex"
$ws.Range("B6").Value = [double]"-1634.544845577873"
$ws.Range("C6").Value = [double]"7.259543146442223e-09"
$ws.Range("A7").Value = "Age"
$ws.Range("B7").Value = [double]"-2.904861095560552"
$ws.Range("C7").Value = [double]"0.8239768622394246"
$ws.Range("A8").Value = "DistSubcenter_res"
$ws.Range("B8").Value = [double]"-463.5664196573799"
$ws.Range("C8").Value = [double]"9.027440263354493e-07"
$ws.Range("A9").Value = "DistCenter_res"
$ws.Range("B9").Value = [double]"790.7879233509702"
$ws.Range("C9").Value = [double]"2.298071237652127e-21"
$ws.Range("A10").Value = "UrbPopDensity_res"
$ws.Range("B10").Value = [double]"0.09013962819305228"
$ws.Range("C10").Value = [double]"0.1252245242826846"
$ws.Range("A11").Value = "UrbBuildDensity_res"
$ws.Range("B11").Value = [double]"-1.27160627876693e-06"
$ws.Range("C11").Value = [double]"0.9588692572482228"
$ws.Range("A12").Value = "IntersecDensity_res"
$ws.Range("B12").Value = [double]"-22.61819907998355"
$ws.Range("C12").Value = [double]"0.1346319978160821"
$ws.Range("A13").Value = "street_length_res"
$ws.Range("B13").Value = [double]"-21.33197352578988"
$ws.Range("C13").Value = [double]"0.1656926581660219"
$ws.Range("A14").Value = "LU_Comm_res"
$ws.Range("B14").Value = [double]"-4158.877505048775"
$ws.Range("C14").Value = [double]"0.05914867602324015"
$ws.Range("A15").Value = "LU_UrbFab_res"
$ws.Range("B15").Value = [double]"-7018.488551213845"
$ws.Range("C15").Value = [double]"0.0004471670249081249"
$ws.Range("A16").Value = "bike_lane_share_res"
$ws.Range("B16").Value = [double]"1548.105043832677"
$ws.Range("C16").Value = [double]"0.6074365041565296"

# --- Sheet 9: summ28015736 -> summ51048321 ---
$ws = $wb.Worksheets.Item(9)
$ws.Name = "summ51048321"

# Drop the old row 3 (Education[T.Secondary+Matura] reference row removed
# from the new model); this shifts all rows below it up by one and
# updates the used range automatically.
$ws.Rows.Item(3).Delete()

# Write the new model's coefficients/p-values (rows 2-16).
$ws.Range("A2").Value = "Intercept"
$ws.Range("B2").Value = [double]"11584.20147062333"
$ws.Range("C2").Value = [double]"6.494413462856286e-06"
$ws.Range("A3").Value = "Education[T.University]"
$ws.Range("B3").Value = [double]"-344.7654445435874"
$ws.Range("C3").Value = [double]"0.310380173240354"
$ws.Range("A4").Value = "Education[T.Unknown/Other]"
$ws.Range("B4").Value = [double]"-330.5727673517993"
$ws.Range("C4").Value = [double]"0.3658471627243087"
$ws.Range("A5").Value = "HHSize"
$ws.Range("B5").Value = [double]"-119.5799685235182"
$ws.Range("C5").Value = [double]"0.3487111491563339"
$ws.Range("A6").Value = "Sex"
$ws.Range("B6").Value = [double]"-1629.09907392681"
$ws.Range("C6").Value = [double]"2.763503507899233e-09"
$ws.Range("A7").Value = "Age"
$ws.Range("B7").Value = [double]"3.866126712177499"
$ws.Range("C7").Value = [double]"0.7575269649104426"
$ws.Range("A8").Value = "DistSubcenter_res"
$ws.Range("B8").Value = [double]"-373.9332256971681"
$ws.Range("C8").Value = [double]"4.451780412788215e-05"
$ws.Range("A9").Value = "DistCenter_res"
$ws.Range("B9").Value = [double]"780.526016590178"
$ws.Range("C9").Value = [double]"1.135987081061133e-22"
$ws.Range("A10").Value = "UrbPopDensity_res"
$ws.Range("B10").Value = [double]"0.08030740285364268"
$ws.Range("C10").Value = [double]"0.1529579239245949"
$ws.Range("A11").Value = "UrbBuildDensity_res"
$ws.Range("B11").Value = [double]"8.358980023061625e-06"
$ws.Range("C11").Value = [double]"0.716526272648601"
$ws.Range("A12").Value = "IntersecDensity_res"
$ws.Range("B12").Value = [double]"-8.237478547613147"
$ws.Range("C12").Value = [double]"0.5780703641704894"
$ws.Range("A13").Value = "street_length_res"
$ws.Range("B13").Value = [double]"-0.4027868557088388"
$ws.Range("C13").Value = [double]"0.9788866964001973"
$ws.Range("A14").Value = "LU_Comm_res"
$ws.Range("B14").Value = [double]"-3571.652215029767"
$ws.Range("C14").Value = [double]"0.09427743884918879"
$ws.Range("A15").Value = "LU_UrbFab_res"
$ws.Range("B15").Value = [double]"-5736.889077964642"
$ws.Range("C15").Value = [double]"0.002959560396638322"
$ws.Range("A16").Value = "bike_lane_share_res"
$ws.Range("B16").Value = [double]"-1977.862540353661"
$ws.Range("C16").Value = [double]"0.4895730972762767"

